$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]'0.000323319659290534'
$ws.Range("C3").Value = [double]'0.0001757085724566559'
$ws.Range("C4").Value = [double]'0.000167343829378347'
$ws.Range("B5").Value = 'type_3'
$ws.Range("C5").Value = [double]'0.0001615497839389279'
$ws.Range("C6").Value = [double]'7.593560334400594e-05'
$ws.Range("B7").Value = 'type_0'
$ws.Range("C7").Value = [double]'7.007430859284367e-05'
$ws.Range("B8").Value = '(h$_{p}$)$_{3}$'
$ws.Range("C8").Value = [double]'5.815195854567683e-05'
$ws.Range("C9").Value = [double]'5.326467721225847e-05'
$ws.Range("C10").Value = [double]'4.929968014736198e-05'
$ws.Range("B11").Value = 'type_1'
$ws.Range("C11").Value = [double]'4.75487959035365e-05'
$ws.Range("C12").Value = [double]'3.921608327345163e-05'
$ws.Range("C13").Value = [double]'2.824189854324345e-05'
$ws.Range("C14").Value = [double]'2.715463632106607e-05'
$ws.Range("C15").Value = [double]'2.455300804959925e-05'
$ws.Range("B16").Value = 'type_2'
$ws.Range("C16").Value = [double]'2.403831982194308e-05'
$ws.Range("C17").Value = [double]'2.369948219689375e-05'
$ws.Range("C18").Value = [double]'2.331067635768819e-05'
$ws.Range("C19").Value = [double]'2.069220219422768e-05'
$ws.Range("C20").Value = [double]'1.978033023138282e-05'
$ws.Range("C21").Value = [double]'1.696587663503903e-05'
$ws.Range("C22").Value = [double]'1.69295312454462e-05'
$ws.Range("B23").Value = '(h$_{p}$)$_{0}$'
$ws.Range("C23").Value = [double]'1.647028120560167e-05'
$ws.Range("A24").Value = 48
$ws.Range("B24").Value = '$(F_{r}^{\text{SCF}})_{1}$'
$ws.Range("C24").Value = [double]'1.4399014374284e-05'
$ws.Range("A25").Value = 91
$ws.Range("B25").Value = '$(\langle pq \vert qp \rangle)_{2}$'
$ws.Range("C25").Value = [double]'1.43312779978816e-05'
$ws.Range("C26").Value = [double]'1.267050212049094e-05'
$ws.Range("B27").Value = '(h$_{p}$)$_{1}$'
$ws.Range("C27").Value = [double]'1.111874518996097e-05'
$ws.Range("A28").Value = 93
$ws.Range("B28").Value = '$(\langle rs \vert sr \rangle)_{2}$'
$ws.Range("C28").Value = [double]'1.075352999252778e-05'
$ws.Range("A29").Value = 58
$ws.Range("B29").Value = '$(F_{r})_{2}$'
$ws.Range("C29").Value = [double]'1.06922904798782e-05'
$ws.Range("C30").Value = [double]'9.97856588650984e-06'
$ws.Range("C31").Value = [double]'9.880168318315773e-06'
$ws.Range("B32").Value = '(h$_{pr}$)$_{3}$'
$ws.Range("C32").Value = [double]'9.412418613818235e-06'
$ws.Range("B33").Value = '(h$_{r}$)$_{1}$'
$ws.Range("C33").Value = [double]'9.213635861829751e-06'
$ws.Range("C34").Value = [double]'8.669854634485175e-06'
$ws.Range("B35").Value = '(h$_{r}$)$_{2}$'
$ws.Range("C35").Value = [double]'7.999564397386185e-06'
$ws.Range("C36").Value = [double]'7.554486870731305e-06'
$ws.Range("A37").Value = 17
$ws.Range("B37").Value = '(h$_{r}$)$_{3}$'
$ws.Range("C37").Value = [double]'7.16988903912837e-06'
$ws.Range("A38").Value = 67
$ws.Range("B38").Value = '$(\eta_{r})_{3}$'
$ws.Range("C38").Value = [double]'7.167021313490361e-06'
$ws.Range("C39").Value = [double]'7.018979581569162e-06'
$ws.Range("C40").Value = [double]'6.997121490467267e-06'
$ws.Range("B41").Value = '(h$_{p}$)$_{2}$'
$ws.Range("C41").Value = [double]'6.874651434815424e-06'
$ws.Range("A42").Value = 13
$ws.Range("B42").Value = 'h$_{qs}$'
$ws.Range("C42").Value = [double]'6.058142791082416e-06'
$ws.Range("A43").Value = 52
$ws.Range("B43").Value = '$(F_{p}^{\text{SCF}})_{2}$'
$ws.Range("C43").Value = [double]'6.038100453020164e-06'
$ws.Range("B44").Value = '(h$_{pr}$)$_{2}$'
$ws.Range("C44").Value = [double]'5.866224091976741e-06'
$ws.Range("C45").Value = [double]'5.584662610670774e-06'
$ws.Range("C46").Value = [double]'5.345206571897236e-06'
$ws.Range("C47").Value = [double]'4.693567335085946e-06'
$ws.Range("A48").Value = 44
$ws.Range("B48").Value = '$(F_{p}^{\text{SCF}})_{1}$'
$ws.Range("C48").Value = [double]'4.297232637943183e-06'
$ws.Range("A49").Value = 66
$ws.Range("B49").Value = '$(F_{r})_{3}$'
$ws.Range("C49").Value = [double]'4.267782444560504e-06'
$ws.Range("A50").Value = 74
$ws.Range("B50").Value = '$(\langle pq \vert pq \rangle)_{0}$'
$ws.Range("C50").Value = [double]'3.967105339601254e-06'
$ws.Range("A51").Value = 36
$ws.Range("B51").Value = '$(F_{r}^{\text{SCF}})_{0}$'
$ws.Range("C51").Value = [double]'3.935547760671783e-06'
$ws.Range("A52").Value = 39
$ws.Range("B52").Value = '$(\eta_{r})_{0}$'
$ws.Range("C52").Value = [double]'3.907724065533214e-06'
$ws.Range("B53").Value = '(h$_{r}$)$_{0}$'
$ws.Range("C53").Value = [double]'3.553340558405851e-06'
$ws.Range("C54").Value = [double]'3.514055150533764e-06'
$ws.Range("C55").Value = [double]'3.43653622366948e-06'
$ws.Range("C56").Value = [double]'3.271165589167211e-06'
$ws.Range("B57").Value = '(h$_{pr}$)$_{0}$'
$ws.Range("C57").Value = [double]'3.0629944695756e-06'
$ws.Range("A58").Value = 9
$ws.Range("B58").Value = '(h$_{pr}$)$_{1}$'
$ws.Range("C58").Value = [double]'2.982622927976264e-06'
$ws.Range("A59").Value = 30
$ws.Range("B59").Value = '$(F_{p})_{0}$'
$ws.Range("C59").Value = [double]'2.912765229541939e-06'
$ws.Range("A60").Value = 70
$ws.Range("B60").Value = '$(\langle pp \vert pp \rangle)_{0}$'
$ws.Range("C60").Value = [double]'2.879258890548974e-06'
$ws.Range("B61").Value = '(h$_{pq}$)$_{1}$'
$ws.Range("C61").Value = [double]'2.857746904189719e-06'
$ws.Range("A62").Value = 60
$ws.Range("B62").Value = '$(F_{p}^{\text{SCF}})_{3}$'
$ws.Range("C62").Value = [double]'2.74927901996683e-06'
$ws.Range("A63").Value = 101
$ws.Range("B63").Value = '$(\langle rs \vert sr \rangle)_{3}$'
$ws.Range("C63").Value = [double]'2.69109803998534e-06'
$ws.Range("C64").Value = [double]'2.601342485196746e-06'
$ws.Range("C65").Value = [double]'2.304631293351454e-06'
$ws.Range("C66").Value = [double]'2.274781968690367e-06'
$ws.Range("C67").Value = [double]'2.161522061566603e-06'
$ws.Range("B68").Value = '(h$_{rs}$)$_{1}$'
$ws.Range("C68").Value = [double]'2.034060694932896e-06'
$ws.Range("C69").Value = [double]'1.867379770044898e-06'
$ws.Range("C70").Value = [double]'1.616492031364258e-06'
$ws.Range("C71").Value = [double]'1.541625136547503e-06'
$ws.Range("C72").Value = [double]'1.420549705372149e-06'
$ws.Range("C73").Value = [double]'1.419893049000004e-06'
$ws.Range("C74").Value = [double]'1.331841715321663e-06'
$ws.Range("B75").Value = '(h$_{rs}$)$_{3}$'
$ws.Range("C75").Value = [double]'1.287323041683143e-06'
$ws.Range("C76").Value = [double]'1.1524700653601e-06'
$ws.Range("B77").Value = '(h$_{rs}$)$_{0}$'
$ws.Range("C77").Value = [double]'9.693119354045891e-07'
$ws.Range("C78").Value = [double]'7.53046982235866e-07'
$ws.Range("B79").Value = '(h$_{pq}$)$_{0}$'
$ws.Range("C79").Value = [double]'6.791581286895688e-07'
$ws.Range("C80").Value = [double]'6.365518206814005e-07'
$ws.Range("B81").Value = '(h$_{pq}$)$_{3}$'
$ws.Range("C81").Value = [double]'5.311686981222513e-07'
$ws.Range("C82").Value = [double]'3.503942292561819e-07'
$ws.Range("B83").Value = '(h$_{rs}$)$_{2}$'
$ws.Range("C83").Value = [double]'3.897830628535145e-08'
$ws.Range("B84").Value = '(h$_{pq}$)$_{2}$'
$ws.Range("C84").Value = [double]'1.084715341690152e-08'
$ws.Range("C85").Value = [double]'3.229006925045185e-09'
$ws.Range("A86").Value = 57
$ws.Range("B86").Value = '$(\omega_{r})_{2}$'
$ws.Range("C86").Value = [double]'2.586757821491536e-09'
$ws.Range("A87").Value = 79
$ws.Range("B87").Value = '$(\langle pq \vert sr \rangle)_{1}$'
$ws.Range("C87").Value = [double]'1.864466745935933e-09'
$ws.Range("C88").Value = [double]'1.782015059132638e-09'
$ws.Range("A89").Value = 49
$ws.Range("B89").Value = '$(\omega_{r})_{1}$'
$ws.Range("C89").Value = [double]'1.75481771446811e-09'
$ws.Range("A90").Value = 69
$ws.Range("B90").Value = '$(\langle pq \vert sr \rangle)_{0}$'
$ws.Range("C90").Value = [double]'1.671646071409022e-09'
$ws.Range("C91").Value = [double]'1.63652090483557e-09'
$ws.Range("C92").Value = [double]'1.450305140952312e-09'
$ws.Range("A93").Value = 33
$ws.Range("B93").Value = '$\omega_{q}$'
$ws.Range("C93").Value = [double]'1.06117104656905e-09'
$ws.Range("A94").Value = 87
$ws.Range("B94").Value = '$(\langle pq \vert sr \rangle)_{2}$'
$ws.Range("C94").Value = [double]'1.049725565464064e-09'
$ws.Range("C95").Value = [double]'2.091820723572891e-10'
